# Fixed few problems with parser.
#
# 1. classNumberOfLines: org.andante.eureka.EurekaApplication's line count
#    was wrong (parser bug) -- correct it from 6 to 3.
# 2. methodNumberOfLines: the parser had missed the implicit no-arg
#    constructor EurekaApplication() -- add its row (0 lines).

$wb = $excel.ActiveWorkbook

$wsClassLines = $wb.Worksheets.Item("classNumberOfLines")
$cell = $wsClassLines.Cells.Item(2, 2)
$cell.NumberFormat = "@"
$cell.Value = "3"

$wsMethodLines = $wb.Worksheets.Item("methodNumberOfLines")
$wsMethodLines.Cells.Item(3, 1).Value = "org.andante.eureka.EurekaApplication"
$wsMethodLines.Cells.Item(3, 2).Value = "EurekaApplication()"
$newLinesCell = $wsMethodLines.Cells.Item(3, 3)
$newLinesCell.NumberFormat = "@"
$newLinesCell.Value = "0"
